$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1): C1, D1, E1 text values
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Update data row (row 2): C2 becomes a string, D2 stays the same string, E2 becomes numeric 1
$ws.Range("C2").Value = "s__CAG-1031 sp000431215"
$ws.Range("D2").Value = "s__CAG-1031 sp000431215"
$ws.Range("E2").Value = 1
